$wb = $excel.ActiveWorkbook

# --- Errors sheet: move selection from B22 to B19 ---
$wsErrors = $wb.Worksheets.Item("Errors")
$wsErrors.Range("B19").Select() | Out-Null
$wsErrors.PageSetup.Orientation = 1

# --- Functional Components sheet: fill in Time Invested (H) / Actual Time (J) hours ---
$wsFunc = $wb.Worksheets.Item("Functional Components")

$hours = @{
  6  = @(0.5, 0.3)
  12 = @(2,   3)
  13 = @(2,   0.8)
  14 = @(1,   1)
  15 = @(2,   3)
  16 = @(0.5, 0.5)
  17 = @(2,   1.5)
  18 = @(0.5, 0.5)
  19 = @(2,   2)
  20 = @(1,   1)
  21 = @(0.5, 0.7)
  22 = @(1,   2)
  23 = @(0.5, 1)
  24 = @(1,   2)
  25 = @(1,   1)
  26 = @(0,   0)
  27 = @(0,   0)
  28 = @(1.5, 2)
  29 = @(0,   0)
  30 = @(0,   0)
  31 = @(0,   0)
}

foreach ($row in $hours.Keys) {
  $vals = $hours[$row]
  $wsFunc.Range("H$row").Value = $vals[0]
  $wsFunc.Range("J$row").Value = $vals[1]
}

$wsFunc.PageSetup.Orientation = 1

# Make "Functional Components" the active/selected sheet
$wsFunc.Activate() | Out-Null

# --- Add Sources sheet: move selection from B7 to B6 (no longer the active tab) ---
$wsAdd = $wb.Worksheets.Item("Add Sources")
$wsAdd.Range("B6").Select() | Out-Null

# Re-activate Functional Components so it remains the active tab after selections above
$wsFunc.Activate() | Out-Null
